# Insert a new data row at row 8 (pushing the existing rows 8..113 down to
# 9..114) and populate it with a new "Achicoria" price-report record for
# Vega Modelo de Temuco, consistent with the other rows in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 8..113 down one row, inserting a new blank row 8.
$ws.Rows.Item(8).Insert()

# Fill the newly inserted row 8 with the new record's data.
$ws.Cells.Item(8, 1).Value = 10
$ws.Cells.Item(8, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(8, 3).Value = "La Araucanía"
$ws.Cells.Item(8, 4).Value = 45083
$ws.Cells.Item(8, 5).Value = 9
$ws.Cells.Item(8, 6).Value = 100112010
$ws.Cells.Item(8, 7).Value = "Achicoria"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 80
$ws.Cells.Item(8, 11).Value = 8000
$ws.Cells.Item(8, 12).Value = 8000
$ws.Cells.Item(8, 13).Value = 8000
$ws.Cells.Item(8, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(8, 15).Value = "Región del Maule"
$ws.Cells.Item(8, 16).Value = 444
$ws.Cells.Item(8, 17).Value = 18
$ws.Cells.Item(8, 18).Value = "Hortaliza"
